# Apply "nu people data partial" update to the raw-data sheet.
#
# Adds previously-missing contact/address details for a couple of people
# (rows 16 & 17 -> phone / e-mail / address / department / city), and
# back-fills a new "neighbourhood" style column (L) for a handful of the
# most recently added rows (28-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (id 66818737) ---------------------------------------------
$ws.Range("E16").Value = "+57 316 4788629"
$ws.Range("F16").Value = "adriabego@yahoo.es"
$ws.Range("G16").Value = "C13-9-36"
$ws.Range("H16").Value = "Valle"
$ws.Range("I16").Value = "Cali"

# --- Row 17 (id 66858587) ---------------------------------------------
$ws.Range("E17").Value = "+57 304 5235781"
$ws.Range("F17").Value = "sicolog@hotmail.com"
$ws.Range("G17").Value = "C4S-48-10-A603"
$ws.Range("H17").Value = "Valle"
$ws.Range("I17").Value = "Cali"

# --- New column L: neighbourhood, filled bottom-up for rows 30,31,29,28
$ws.Range("L30").Value = "Villa del Rosario"
$ws.Range("L31").Value = "Villa del Rosario"
$ws.Range("L29").Value = "Santa Ana"
$ws.Range("L28").Value = "Metropolitano Norte"

# Match the workbook's saved selection state.
$ws.Range("G16:G17").Select()
